$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: cardholder name and card number
$ws.Range("C2").Value = "Hartmut"
# Leading apostrophe forces this purely-numeric-looking string to be stored
# as text (matching the source workbook, where the card number is text)
$ws.Range("B3").Value = "'2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# Opening balance line
$ws.Range("D5").Value = "KONTOSTAND AM 27.11.2024"

# Row 6
$ws.Range("B6").Value = "28.11."
$ws.Range("C6").Value = "29.11."
$ws.Range("D6").Value = "BURGER KING Wernigerode"
$ws.Range("E6").Value = "22,39-"

# Row 7
$ws.Range("B7").Value = "02.12."
$ws.Range("C7").Value = "03.12."
$ws.Range("D7").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 11911497"
$ws.Range("E7").Value = "85,63-"

# Row 8
$ws.Range("B8").Value = "06.12."
$ws.Range("C8").Value = "07.12."
$ws.Range("D8").Value = "BEITRAG Allianz SE K-99275379"
$ws.Range("E8").Value = "54,66-"

# Row 9 - cleared out (previously had a transaction, now blank)
$ws.Range("B9").Value = ""
$ws.Range("C9").Value = ""
$ws.Range("D9").Value = ""
$ws.Range("E9").Value = ""
$ws.Range("E9").HorizontalAlignment = -4108
$ws.Range("E9").VerticalAlignment = -4108
$ws.Range("E9").WrapText = $true

# Closing balance line
$ws.Range("D12").Value = "KONTOSTAND AM 09.12.2024"
$ws.Range("E12").Value = "162,68-"

# Next statement date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 14.12.2024"
